$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.019276059405407
$ws.Cells.Item(2, 4).Value = 1.025259596646197
$ws.Cells.Item(2, 5).Value = 1.022891904291117
$ws.Cells.Item(2, 6).Value = 1.030488057141025
$ws.Cells.Item(2, 9).Value = 1.02902564449912
$ws.Cells.Item(2, 10).Value = 1.024479743265611
$ws.Cells.Item(2, 11).Value = 1.028085690362026
$ws.Cells.Item(2, 12).Value = 1.025724949130821
$ws.Cells.Item(2, 13).Value = 1.033298924045309
$ws.Cells.Item(2, 14).Value = 1.012222202112561

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.020219859434824
$ws.Cells.Item(3, 4).Value = 1.025951695828123
$ws.Cells.Item(3, 5).Value = 1.023782517279035
$ws.Cells.Item(3, 6).Value = 1.031650389478844
$ws.Cells.Item(3, 9).Value = 1.029181168928773
$ws.Cells.Item(3, 10).Value = 1.025060257866242
$ws.Cells.Item(3, 11).Value = 1.028585367592407
$ws.Cells.Item(3, 12).Value = 1.026422095978817
$ws.Cells.Item(3, 13).Value = 1.034268671050994
$ws.Cells.Item(3, 14).Value = 1.012416974006674

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.020830917333048
$ws.Cells.Item(4, 4).Value = 1.026399602648369
$ws.Cells.Item(4, 5).Value = 1.024359525066087
$ws.Cells.Item(4, 6).Value = 1.032402990707849
$ws.Cells.Item(4, 9).Value = 1.029280482458022
$ws.Cells.Item(4, 10).Value = 1.025435662004712
$ws.Cells.Item(4, 11).Value = 1.028908090695391
$ws.Cells.Item(4, 12).Value = 1.026873288201085
$ws.Cells.Item(4, 13).Value = 1.0348960860374
$ws.Cells.Item(4, 14).Value = 1.012542852625679

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.021087890722084
$ws.Cells.Item(5, 4).Value = 1.026587918828151
$ws.Cells.Item(5, 5).Value = 1.024602270895171
$ws.Cells.Item(5, 6).Value = 1.032719502701747
$ws.Cells.Item(5, 9).Value = 1.029321917162517
$ws.Cells.Item(5, 10).Value = 1.025593426732428
$ws.Cells.Item(5, 11).Value = 1.029043618694319
$ws.Cells.Item(5, 12).Value = 1.027062990531525
$ws.Cells.Item(5, 13).Value = 1.03515983296484
$ws.Cells.Item(5, 14).Value = 1.012595735358397

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.021131042630766
$ws.Cells.Item(6, 4).Value = 1.026619538856419
$ws.Cells.Item(6, 5).Value = 1.024643039030648
$ws.Cells.Item(6, 6).Value = 1.032772653405233
$ws.Cells.Item(6, 9).Value = 1.029328855648165
$ws.Cells.Item(6, 10).Value = 1.025619912869395
$ws.Cells.Item(6, 11).Value = 1.0290663659277
$ws.Cells.Item(6, 12).Value = 1.027094843602567
$ws.Cells.Item(6, 13).Value = 1.03520411615672
$ws.Cells.Item(6, 14).Value = 1.012604612444893

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.020834350692204
$ws.Cells.Item(7, 4).Value = 1.02640211887639
$ws.Cells.Item(7, 5).Value = 1.024362767974113
$ws.Cells.Item(7, 6).Value = 1.032407219491148
$ws.Cells.Item(7, 9).Value = 1.029281037355511
$ws.Cells.Item(7, 10).Value = 1.025437770282364
$ws.Cells.Item(7, 11).Value = 1.02890990219709
$ws.Cells.Item(7, 12).Value = 1.026875822930934
$ws.Cells.Item(7, 13).Value = 1.034899610308907
$ws.Cells.Item(7, 14).Value = 1.012543559391674

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.019594946924092
$ws.Cells.Item(8, 4).Value = 1.02549347911589
$ws.Cells.Item(8, 5).Value = 1.023192741098068
$ws.Cells.Item(8, 6).Value = 1.030880770443653
$ws.Cells.Item(8, 9).Value = 1.029078478019425
$ws.Cells.Item(8, 10).Value = 1.024675977616756
$ws.Cells.Item(8, 11).Value = 1.028254682681544
$ws.Cells.Item(8, 12).Value = 1.02596053372121
$ws.Cells.Item(8, 13).Value = 1.033626670298701
$ws.Cells.Item(8, 14).Value = 1.012288057445782

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.017413714123701
$ws.Cells.Item(9, 4).Value = 1.023892945643361
$ws.Cells.Item(9, 5).Value = 1.021136575270754
$ws.Cells.Item(9, 6).Value = 1.028194756430287
$ws.Cells.Item(9, 9).Value = 1.0287114420418
$ws.Cells.Item(9, 10).Value = 1.023331888046953
$ws.Cells.Item(9, 11).Value = 1.027095530691469
$ws.Cells.Item(9, 12).Value = 1.024348414208576
$ws.Cells.Item(9, 13).Value = 1.031383007926253
$ws.Cells.Item(9, 14).Value = 1.011836680462803

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.015961445988429
$ws.Cells.Item(10, 4).Value = 1.02282639931614
$ws.Cells.Item(10, 5).Value = 1.019769601959864
$ws.Cells.Item(10, 6).Value = 1.026406617579014
$ws.Cells.Item(10, 9).Value = 1.028459984225308
$ws.Cells.Item(10, 10).Value = 1.022434715267816
$ws.Cells.Item(10, 11).Value = 1.026319735482135
$ws.Cells.Item(10, 12).Value = 1.023274210076371
$ws.Cells.Item(10, 13).Value = 1.029886846025275
$ws.Cells.Item(10, 14).Value = 1.011535006303833

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015333052159977
$ws.Cells.Item(11, 4).Value = 1.022364700047922
$ws.Cells.Item(11, 5).Value = 1.019178601527913
$ws.Cells.Item(11, 6).Value = 1.025632932930567
$ws.Cells.Item(11, 9).Value = 1.028349499984653
$ws.Cells.Item(11, 10).Value = 1.022045973410058
$ws.Cells.Item(11, 11).Value = 1.025983099047923
$ws.Cells.Item(11, 12).Value = 1.022809205886587
$ws.Cells.Item(11, 13).Value = 1.029238899263672
$ws.Cells.Item(11, 14).Value = 1.011404202166334

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015099706464956
$ws.Cells.Item(12, 4).Value = 1.022193223866915
$ws.Cells.Item(12, 5).Value = 1.018959214985286
$ws.Cells.Item(12, 6).Value = 1.025345640502232
$ws.Cells.Item(12, 9).Value = 1.028308220987426
$ws.Cells.Item(12, 10).Value = 1.021901538879983
$ws.Cells.Item(12, 11).Value = 1.025857951037411
$ws.Cells.Item(12, 12).Value = 1.022636503364842
$ws.Cells.Item(12, 13).Value = 1.028998208003185
$ws.Cells.Item(12, 14).Value = 1.01135558932024

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015149756835197
$ws.Cells.Item(13, 4).Value = 1.022230005190395
$ws.Cells.Item(13, 5).Value = 1.019006267914219
$ws.Cells.Item(13, 6).Value = 1.025407261703921
$ws.Cells.Item(13, 9).Value = 1.028317086337732
$ws.Cells.Item(13, 10).Value = 1.021932522309575
$ws.Cells.Item(13, 11).Value = 1.025884800507008
$ws.Cells.Item(13, 12).Value = 1.022673547692081
$ws.Cells.Item(13, 13).Value = 1.029049837775876
$ws.Cells.Item(13, 14).Value = 1.011366018129629

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.0153137623472
$ws.Cells.Item(14, 4).Value = 1.022350525369448
$ws.Cells.Item(14, 5).Value = 1.019160464161063
$ws.Cells.Item(14, 6).Value = 1.025609183445392
$ws.Cells.Item(14, 9).Value = 1.028346092748457
$ws.Cells.Item(14, 10).Value = 1.022034035188654
$ws.Cells.Item(14, 11).Value = 1.025972756437653
$ws.Cells.Item(14, 12).Value = 1.022794929815075
$ws.Cells.Item(14, 13).Value = 1.029219003955459
$ws.Cells.Item(14, 14).Value = 1.011400184349641

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015414820547843
$ws.Cells.Item(15, 4).Value = 1.022424784444329
$ws.Cells.Item(15, 5).Value = 1.019255487782936
$ws.Cells.Item(15, 6).Value = 1.025733605811568
$ws.Cells.Item(15, 9).Value = 1.02836393273309
$ws.Cells.Item(15, 10).Value = 1.022096575540307
$ws.Cells.Item(15, 11).Value = 1.026026934931739
$ws.Cells.Item(15, 12).Value = 1.022869720116592
$ws.Cells.Item(15, 13).Value = 1.029323230832267
$ws.Cells.Item(15, 14).Value = 1.01142123179959

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016003159892845
$ws.Cells.Item(16, 4).Value = 1.022857043447787
$ws.Cells.Item(16, 5).Value = 1.019808843921033
$ws.Cells.Item(16, 6).Value = 1.026457976909924
$ws.Cells.Item(16, 9).Value = 1.028467282992092
$ws.Cells.Item(16, 10).Value = 1.022460509354065
$ws.Cells.Item(16, 11).Value = 1.026342061994115
$ws.Cells.Item(16, 12).Value = 1.023305073742128
$ws.Cells.Item(16, 13).Value = 1.029929846024376
$ws.Cells.Item(16, 14).Value = 1.011543683631889

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.016372329852826
$ws.Cells.Item(17, 4).Value = 1.023128221633079
$ws.Cells.Item(17, 5).Value = 1.02015619335583
$ws.Cells.Item(17, 6).Value = 1.026912513826315
$ws.Cells.Item(17, 9).Value = 1.028531683402293
$ws.Cells.Item(17, 10).Value = 1.022688726191046
$ws.Cells.Item(17, 11).Value = 1.026539542762381
$ws.Cells.Item(17, 12).Value = 1.023578195758333
$ws.Cells.Item(17, 13).Value = 1.030310333158959
$ws.Cells.Item(17, 14).Value = 1.011620447165407

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.016587703624927
$ws.Cells.Item(18, 4).Value = 1.023286407009615
$ws.Cells.Item(18, 5).Value = 1.020358883955358
$ws.Cells.Item(18, 6).Value = 1.027177694591008
$ws.Cells.Item(18, 9).Value = 1.028569092541519
$ws.Cells.Item(18, 10).Value = 1.022821816073483
$ws.Cells.Item(18, 11).Value = 1.02665466120865
$ws.Cells.Item(18, 12).Value = 1.023737516022037
$ws.Cells.Item(18, 13).Value = 1.030532255345188
$ws.Cells.Item(18, 14).Value = 1.011665204917512

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.016661147825439
$ws.Cells.Item(19, 4).Value = 1.023340346090904
$ws.Cells.Item(19, 5).Value = 1.020428011034252
$ws.Cells.Item(19, 6).Value = 1.027268124069425
$ws.Cells.Item(19, 9).Value = 1.028581821869469
$ws.Cells.Item(19, 10).Value = 1.022867192002491
$ws.Cells.Item(19, 11).Value = 1.026693901941738
$ws.Cells.Item(19, 12).Value = 1.02379184226627
$ws.Cells.Item(19, 13).Value = 1.030607923464985
$ws.Cells.Item(19, 14).Value = 1.011680463242793

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.01633271694013
$ws.Cells.Item(20, 4).Value = 1.02309912555899
$ws.Cells.Item(20, 5).Value = 1.020118916986334
$ws.Cells.Item(20, 6).Value = 1.026863740415978
$ws.Cells.Item(20, 9).Value = 1.028524789829
$ws.Cells.Item(20, 10).Value = 1.022664243280074
$ws.Cells.Item(20, 11).Value = 1.026518362050252
$ws.Cells.Item(20, 12).Value = 1.023548891020337
$ws.Cells.Item(20, 13).Value = 1.030269511476128
$ws.Cells.Item(20, 14).Value = 1.011612212933457

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.01526546493759
$ws.Cells.Item(21, 4).Value = 1.022315034633871
$ws.Cells.Item(21, 5).Value = 1.019115053411256
$ws.Cells.Item(21, 6).Value = 1.025549720084964
$ws.Cells.Item(21, 9).Value = 1.028337557712172
$ws.Cells.Item(21, 10).Value = 1.022004143233865
$ws.Cells.Item(21, 11).Value = 1.02594685853558
$ws.Cells.Item(21, 12).Value = 1.022759185227002
$ws.Cells.Item(21, 13).Value = 1.029169189146663
$ws.Cells.Item(21, 14).Value = 1.011390123975126

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014594833470894
$ws.Cells.Item(22, 4).Value = 1.02182215882486
$ws.Cells.Item(22, 5).Value = 1.018484679550397
$ws.Cells.Item(22, 6).Value = 1.024724055212224
$ws.Cells.Item(22, 9).Value = 1.028218447722486
$ws.Cells.Item(22, 10).Value = 1.021588889703179
$ws.Cells.Item(22, 11).Value = 1.025586917217735
$ws.Cells.Item(22, 12).Value = 1.022262786192345
$ws.Cells.Item(22, 13).Value = 1.028477286069754
$ws.Cells.Item(22, 14).Value = 1.011250335419997

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.014950310564912
$ws.Cells.Item(23, 4).Value = 1.022083430510932
$ws.Cells.Item(23, 5).Value = 1.018818776999096
$ws.Cells.Item(23, 6).Value = 1.025161707375025
$ws.Cells.Item(23, 9).Value = 1.028281721793501
$ws.Cells.Item(23, 10).Value = 1.021809044345156
$ws.Cells.Item(23, 11).Value = 1.025777786916678
$ws.Cells.Item(23, 12).Value = 1.022525925067523
$ws.Cells.Item(23, 13).Value = 1.028844085252063
$ws.Cells.Item(23, 14).Value = 1.011324454337407

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016350616176806
$ws.Cells.Item(24, 4).Value = 1.023112272786295
$ws.Cells.Item(24, 5).Value = 1.0201357603027
$ws.Cells.Item(24, 6).Value = 1.026885778844123
$ws.Cells.Item(24, 9).Value = 1.028527905215528
$ws.Cells.Item(24, 10).Value = 1.022675306131373
$ws.Cells.Item(24, 11).Value = 1.026527932914798
$ws.Cells.Item(24, 12).Value = 1.023562132530469
$ws.Cells.Item(24, 13).Value = 1.0302879570667
$ws.Cells.Item(24, 14).Value = 1.011615933681345

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.017977285034408
$ws.Cells.Item(25, 4).Value = 1.024306643029241
$ws.Cells.Item(25, 5).Value = 1.021667477224184
$ws.Cells.Item(25, 6).Value = 1.028888707735012
$ws.Cells.Item(25, 9).Value = 1.028807524148931
$ws.Cells.Item(25, 10).Value = 1.023679566461761
$ws.Cells.Item(25, 11).Value = 1.027395736226568
$ws.Cells.Item(25, 12).Value = 1.024765093594183
$ws.Cells.Item(25, 13).Value = 1.031963116839774
$ws.Cells.Item(25, 14).Value = 1.011953506752688
